# "added paged numbers to pinyin"
# Prepend a leading space marker to the chinese / pinyin / zhuyin cells for
# the five "向（諸位）點傳師…駕" rows (rows 27, 29, 31, 33, 35). The Chinese
# column gets a full-width space (U+3000) and the pinyin / zhuyin columns
# get a regular ASCII space, making room for page numbers to be filled in
# later.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(27, 29, 31, 33, 35)

foreach ($r in $rows) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $cCell.Value = [string]::Concat([char]0x3000, $cCell.Value2)
    $dCell.Value = " " + $dCell.Value2
    $eCell.Value = " " + $eCell.Value2
}

$ws.Range("M28").Select() | Out-Null
